# Update region_names worksheet with today's numbers:
# Insert four new country rows (keeping the existing alphabetical ordering
# of column A) with their corresponding region in column B.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("countries")

# Each entry: (row at which to insert a new row, country name, region)
# Rows are given in ascending final-position order so that inserting from
# top to bottom (each insertion shifting everything below it down by one)
# lands every subsequent entry on the correct row.
$newEntries = @(
    @{ Row = 26;  Name = "Bonaire, Saint Eustatius and Saba"; Region = "north america" },
    @{ Row = 119; Name = "Malawi";                             Region = "africa" },
    @{ Row = 165; Name = "Sao_Tome_and_Principe";              Region = "africa" },
    @{ Row = 205; Name = "Yemen";                               Region = "asia" }
)

foreach ($entry in $newEntries) {
    $r = $entry.Row
    $ws.Rows.Item($r).Insert()
    $ws.Cells.Item($r, 1).Value = $entry.Name
    $ws.Cells.Item($r, 2).Value = $entry.Region
}

# Restore the scroll position / selection recorded by Excel after the edit.
$excel.Goto($ws.Range("A189"))
$ws.Range("A205").Select()
